$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: record the "out" (check-out) timestamp and flip status to OUT
$ws.Range("E2").Value = "2023-07-05 15:43:23"
$ws.Range("F2").Value = "OUT"
